$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("fumarate_mvacid")
$ws1.Range('F2').Value = 'RXN-863_reverse: [O87942, O87943, O87944, Q8VPT6, Q8VPT8] | R05598: [] | 4.1.99.11_13: [] | 4.1.99.11_07: [] | RXN-942_reverse: [Q8VPT6, Q8VPT7] | 4.1.99.11_08: [] | '
$ws1.Range('G2').Value = 'R02243: [] | ACONITATE-DECARBOXYLASE-RXN_reverse: [B3IUN8, P54987, Q0C8L3] | '
$ws1.Range('H2').Value = '3.7.1.5_0: [] | R01335_reverse: [] | 3.7.1.2_08: [] | RXN-12096_reverse: [] | 3.7.1.2_05: [] | R00324: [] | ACETYLPYRUVATE-HYDROLASE-RXN_reverse: [] | 3.7.1.2_03: [] | '
$ws1.Range('I2').Value = 'RXN-16543_reverse: [] | '
$ws1.Range('F3').Value = 'RXN-863_reverse: [O87942, O87943, O87944, Q8VPT6, Q8VPT8] | R05598: [] | 4.1.99.11_13: [] | 4.1.99.11_07: [] | RXN-942_reverse: [Q8VPT6, Q8VPT7] | 4.1.99.11_08: [] | '
$ws1.Range('G3').Value = '4.1.1.68-RXN_reverse: [] | '
$ws1.Range('H3').Value = 'R02243: [] | ACONITATE-DECARBOXYLASE-RXN_reverse: [B3IUN8, P54987, Q0C8L3] | '
$ws1.Range('I3').Value = 'RXN-16543_reverse: [] | '
$ws1.Range('F4').Value = 'R02243: [] | ACONITATE-DECARBOXYLASE-RXN_reverse: [B3IUN8, P54987, Q0C8L3] | '
$ws1.Range('G4').Value = 'RXN-863_reverse: [O87942, O87943, O87944, Q8VPT6, Q8VPT8] | R05598: [] | 4.1.99.11_13: [] | '
$ws1.Range('H4').Value = '1.14.13.25_24: [] | '
$ws1.Range('I4').Value = '1.1.3.13_53: [] | '
$ws1.Range('F5').Value = 'R02243: [] | ACONITATE-DECARBOXYLASE-RXN_reverse: [B3IUN8, P54987, Q0C8L3] | '
$ws1.Range('G5').Value = 'RXN-863_reverse: [O87942, O87943, O87944, Q8VPT6, Q8VPT8] | R05598: [] | 4.1.99.11_13: [] | '
$ws1.Range('H5').Value = '1.1.3.13_53: [] | '
$ws1.Range('I5').Value = '1.14.13.25_24: [] | '
$ws1.Range('F6').Value = 'RXN0-1382: [A1ADQ1, A4YXN2, A5EGD7, A7ZPI2, A8A2M8, B1IX88, B1LMH0, B1X9P6, B2TWX3, B3QBS6, B5YYX4, B6I6S5, B6JE29, B7LBS7, B7M6P3, B7MH34, B7MY33, B7N5X4, B7NPQ8, B7UG84, C4ZVR1, O06644, O87838, P69902, P69903, Q07Q82, Q0TF87, Q139H7, Q1R8Z2, Q217M3, Q2IUI7, Q31Y97, Q32DG9, Q3YZF6, Q6N8F8, Q82M40, Q89QH2, Q8FFE8, Q8XBR7] | R07290: [] | 2.8.3.16_1: [] | 2.8.3.9_14: [] | '
$ws1.Range('G6').Value = 'RXN-863_reverse: [O87942, O87943, O87944, Q8VPT6, Q8VPT8] | R05598: [] | 4.1.99.11_13: [] | '
$ws1.Range('H6').Value = '1.14.13.25_24: [] | '
$ws1.Range('I6').Value = '1.1.1.71_096: [] | 1.1.1.2_197: [] | 1.1.1.283_27: [] | 1.1.1.184_198: [] | 1.1.1.71_097: [] | 1.1.1.2_199: [] | 1.1.1.184_200: [] | 1.1.1.71_104: [] | 1.1.1.2_211: [] | 1.1.1.283_33: [] | '
$ws1.Range('F7').Value = 'RXN-863_reverse: [O87942, O87943, O87944, Q8VPT6, Q8VPT8] | R05598: [] | 4.1.99.11_13: [] | 4.1.99.11_07: [] | RXN-942_reverse: [Q8VPT6, Q8VPT7] | 4.1.99.11_08: [] | '
$ws1.Range('G7').Value = 'RXN-16948_reverse: [] | '
$ws1.Range('H7').Value = 'RXN0-1382: [A1ADQ1, A4YXN2, A5EGD7, A7ZPI2, A8A2M8, B1IX88, B1LMH0, B1X9P6, B2TWX3, B3QBS6, B5YYX4, B6I6S5, B6JE29, B7LBS7, B7M6P3, B7MH34, B7MY33, B7N5X4, B7NPQ8, B7UG84, C4ZVR1, O06644, O87838, P69902, P69903, Q07Q82, Q0TF87, Q139H7, Q1R8Z2, Q217M3, Q2IUI7, Q31Y97, Q32DG9, Q3YZF6, Q6N8F8, Q82M40, Q89QH2, Q8FFE8, Q8XBR7] | R07290: [] | 2.8.3.16_1: [] | 2.8.3.9_14: [] | '
$ws1.Range('I7').Value = 'RXN-16543_reverse: [] | '
$ws1.Range('F8').Value = 'RXN-863_reverse: [O87942, O87943, O87944, Q8VPT6, Q8VPT8] | R05598: [] | 4.1.99.11_13: [] | 4.1.99.11_07: [] | RXN-942_reverse: [Q8VPT6, Q8VPT7] | 4.1.99.11_08: [] | '
$ws1.Range('G8').Value = '4.1.1.68-RXN_reverse: [] | '
$ws1.Range('H8').Value = 'RXN0-1382: [A1ADQ1, A4YXN2, A5EGD7, A7ZPI2, A8A2M8, B1IX88, B1LMH0, B1X9P6, B2TWX3, B3QBS6, B5YYX4, B6I6S5, B6JE29, B7LBS7, B7M6P3, B7MH34, B7MY33, B7N5X4, B7NPQ8, B7UG84, C4ZVR1, O06644, O87838, P69902, P69903, Q07Q82, Q0TF87, Q139H7, Q1R8Z2, Q217M3, Q2IUI7, Q31Y97, Q32DG9, Q3YZF6, Q6N8F8, Q82M40, Q89QH2, Q8FFE8, Q8XBR7] | R07290: [] | 2.8.3.16_1: [] | 2.8.3.9_14: [] | '
$ws1.Range('I8').Value = 'RXN-16543_reverse: [] | '
$ws1.Range('F9').Value = 'R02243: [] | ACONITATE-DECARBOXYLASE-RXN_reverse: [B3IUN8, P54987, Q0C8L3] | '
$ws1.Range('G9').Value = 'RXN-863_reverse: [O87942, O87943, O87944, Q8VPT6, Q8VPT8] | R05598: [] | 4.1.99.11_13: [] | '
$ws1.Range('H9').Value = 'RXN-8094_reverse: [B0LB00, B0LB01, C4NYZ3, G3X982, H9TB17, H9TB18, H9TB19, O54754, P48034, P80456, Q01772, Q06278, Q3HRQ2, Q3TYQ9, Q5FB27, Q5QE78, Q5QE79, Q5QE80, Q5SGK3, Q69R21, Q7G192, Q7G193, Q7XH05, Q852M1, Q852M2, Q960A1, Q9FV24, Q9FV25, Q9FYG4, Q9Z0U5] | 1.2.3.1_110: [] | 1.2.3.1_111: [] | 1.2.3.1_090: [] | '
$ws1.Range('I9').Value = '1.1.3.13_53: [] | '
$ws1.Range('F10').Value = 'RXN-863_reverse: [O87942, O87943, O87944, Q8VPT6, Q8VPT8] | R05598: [] | 4.1.99.11_13: [] | 4.1.99.11_07: [] | RXN-942_reverse: [Q8VPT6, Q8VPT7] | 4.1.99.11_08: [] | '
$ws1.Range('G10').Value = 'RXN-11727: [] | R09545_reverse: [] | 1.1.1.308_0: [] | '
$ws1.Range('H10').Value = 'RXN0-1382: [A1ADQ1, A4YXN2, A5EGD7, A7ZPI2, A8A2M8, B1IX88, B1LMH0, B1X9P6, B2TWX3, B3QBS6, B5YYX4, B6I6S5, B6JE29, B7LBS7, B7M6P3, B7MH34, B7MY33, B7N5X4, B7NPQ8, B7UG84, C4ZVR1, O06644, O87838, P69902, P69903, Q07Q82, Q0TF87, Q139H7, Q1R8Z2, Q217M3, Q2IUI7, Q31Y97, Q32DG9, Q3YZF6, Q6N8F8, Q82M40, Q89QH2, Q8FFE8, Q8XBR7] | R07290: [] | 2.8.3.16_1: [] | 2.8.3.9_14: [] | '
$ws1.Range('I10').Value = '3.7.1.5_0: [] | R01335_reverse: [] | 3.7.1.2_08: [] | RXN-12096_reverse: [] | 3.7.1.2_05: [] | R00324: [] | ACETYLPYRUVATE-HYDROLASE-RXN_reverse: [] | 3.7.1.2_03: [] | '
$ws1.Range('F11').Value = '1.4.3.21_04: [] | R01469_reverse: [] | ETHANOLAMINE-OXIDASE-RXN: [D5DET5] | 1.4.3.8_2: [] | 1.4.3.4_28: [] | 1.4.3.21_45: [] | RXN-6381: [P19801, P36633, Q12556, Q59118, Q8JZQ5] | R03139_reverse: [] | RXN-12732: [] | RXN-9757: [] | 1.4.3.22_21: [] | 1.4.3.22_42: [] | 1.4.3.21_57: [] | 1.4.3.21_41: [] | 1.4.3.21_38: [] | 1.4.3.8_3: [] | R02529_reverse: [] | AMACETOXID-RXN: [O08590, O23349, O42890, O46406, O70423, O75106, P12807, P46881, P46883, P49250, P49252, P80695, Q07121, Q07123, Q16853, Q29437, Q43077, Q812C9, Q8H1H9, Q9P7F2, Q9TTK6] | 1.4.3.21_50: [] | '
$ws1.Range('G11').Value = 'RXN0-1382: [A1ADQ1, A4YXN2, A5EGD7, A7ZPI2, A8A2M8, B1IX88, B1LMH0, B1X9P6, B2TWX3, B3QBS6, B5YYX4, B6I6S5, B6JE29, B7LBS7, B7M6P3, B7MH34, B7MY33, B7N5X4, B7NPQ8, B7UG84, C4ZVR1, O06644, O87838, P69902, P69903, Q07Q82, Q0TF87, Q139H7, Q1R8Z2, Q217M3, Q2IUI7, Q31Y97, Q32DG9, Q3YZF6, Q6N8F8, Q82M40, Q89QH2, Q8FFE8, Q8XBR7] | R07290: [] | 2.8.3.16_1: [] | 2.8.3.9_14: [] | '
$ws1.Range('H11').Value = 'RXN-863_reverse: [O87942, O87943, O87944, Q8VPT6, Q8VPT8] | R05598: [] | 4.1.99.11_13: [] | '
$ws1.Range('I11').Value = 'RXN-16543_reverse: [] | '
$ws1.Range('F12').Value = '1.4.3.21_04: [] | R01469_reverse: [] | ETHANOLAMINE-OXIDASE-RXN: [D5DET5] | 1.4.3.8_2: [] | 1.4.3.4_28: [] | 1.4.3.21_45: [] | RXN-6381: [P19801, P36633, Q12556, Q59118, Q8JZQ5] | R03139_reverse: [] | RXN-12732: [] | RXN-9757: [] | 1.4.3.22_21: [] | 1.4.3.22_42: [] | 1.4.3.21_57: [] | 1.4.3.21_41: [] | 1.4.3.21_38: [] | 1.4.3.8_3: [] | R02529_reverse: [] | AMACETOXID-RXN: [O08590, O23349, O42890, O46406, O70423, O75106, P12807, P46881, P46883, P49250, P49252, P80695, Q07121, Q07123, Q16853, Q29437, Q43077, Q812C9, Q8H1H9, Q9P7F2, Q9TTK6] | 1.4.3.21_50: [] | '
$ws1.Range('G12').Value = 'RXN0-1382: [A1ADQ1, A4YXN2, A5EGD7, A7ZPI2, A8A2M8, B1IX88, B1LMH0, B1X9P6, B2TWX3, B3QBS6, B5YYX4, B6I6S5, B6JE29, B7LBS7, B7M6P3, B7MH34, B7MY33, B7N5X4, B7NPQ8, B7UG84, C4ZVR1, O06644, O87838, P69902, P69903, Q07Q82, Q0TF87, Q139H7, Q1R8Z2, Q217M3, Q2IUI7, Q31Y97, Q32DG9, Q3YZF6, Q6N8F8, Q82M40, Q89QH2, Q8FFE8, Q8XBR7] | R07290: [] | 2.8.3.16_1: [] | 2.8.3.9_14: [] | '
$ws1.Range('H12').Value = 'RXN-863_reverse: [O87942, O87943, O87944, Q8VPT6, Q8VPT8] | R05598: [] | 4.1.99.11_13: [] | '
$ws1.Range('I12').Value = 'RXN-16543_reverse: [] | '
$ws1.Range('F13').Value = '1.4.3.21_04: [] | R01469_reverse: [] | ETHANOLAMINE-OXIDASE-RXN: [D5DET5] | 1.4.3.8_2: [] | 1.4.3.4_28: [] | 1.4.3.21_45: [] | RXN-6381: [P19801, P36633, Q12556, Q59118, Q8JZQ5] | R03139_reverse: [] | RXN-12732: [] | RXN-9757: [] | 1.4.3.22_21: [] | 1.4.3.22_42: [] | 1.4.3.21_57: [] | 1.4.3.21_41: [] | 1.4.3.21_38: [] | 1.4.3.8_3: [] | R02529_reverse: [] | AMACETOXID-RXN: [O08590, O23349, O42890, O46406, O70423, O75106, P12807, P46881, P46883, P49250, P49252, P80695, Q07121, Q07123, Q16853, Q29437, Q43077, Q812C9, Q8H1H9, Q9P7F2, Q9TTK6] | 1.4.3.21_50: [] | '
$ws1.Range('G13').Value = 'RXN0-1382: [A1ADQ1, A4YXN2, A5EGD7, A7ZPI2, A8A2M8, B1IX88, B1LMH0, B1X9P6, B2TWX3, B3QBS6, B5YYX4, B6I6S5, B6JE29, B7LBS7, B7M6P3, B7MH34, B7MY33, B7N5X4, B7NPQ8, B7UG84, C4ZVR1, O06644, O87838, P69902, P69903, Q07Q82, Q0TF87, Q139H7, Q1R8Z2, Q217M3, Q2IUI7, Q31Y97, Q32DG9, Q3YZF6, Q6N8F8, Q82M40, Q89QH2, Q8FFE8, Q8XBR7] | R07290: [] | 2.8.3.16_1: [] | 2.8.3.9_14: [] | '
$ws1.Range('H13').Value = 'RXN-863_reverse: [O87942, O87943, O87944, Q8VPT6, Q8VPT8] | R05598: [] | 4.1.99.11_13: [] | '
$ws1.Range('I13').Value = 'RXN-16543_reverse: [] | '
$ws1.Range('F14').Value = 'RXN-863_reverse: [O87942, O87943, O87944, Q8VPT6, Q8VPT8] | R05598: [] | 4.1.99.11_13: [] | 4.1.99.11_07: [] | RXN-942_reverse: [Q8VPT6, Q8VPT7] | 4.1.99.11_08: [] | '
$ws1.Range('G14').Value = 'RXN-11727: [] | R09545_reverse: [] | 1.1.1.308_0: [] | '
$ws1.Range('H14').Value = 'RXN0-1382: [A1ADQ1, A4YXN2, A5EGD7, A7ZPI2, A8A2M8, B1IX88, B1LMH0, B1X9P6, B2TWX3, B3QBS6, B5YYX4, B6I6S5, B6JE29, B7LBS7, B7M6P3, B7MH34, B7MY33, B7N5X4, B7NPQ8, B7UG84, C4ZVR1, O06644, O87838, P69902, P69903, Q07Q82, Q0TF87, Q139H7, Q1R8Z2, Q217M3, Q2IUI7, Q31Y97, Q32DG9, Q3YZF6, Q6N8F8, Q82M40, Q89QH2, Q8FFE8, Q8XBR7] | R07290: [] | 2.8.3.16_1: [] | 2.8.3.9_14: [] | '
$ws1.Range('I14').Value = '4.1.1.68-RXN_reverse: [] | '
$ws1.Range('F15').Value = 'RXN0-1382: [A1ADQ1, A4YXN2, A5EGD7, A7ZPI2, A8A2M8, B1IX88, B1LMH0, B1X9P6, B2TWX3, B3QBS6, B5YYX4, B6I6S5, B6JE29, B7LBS7, B7M6P3, B7MH34, B7MY33, B7N5X4, B7NPQ8, B7UG84, C4ZVR1, O06644, O87838, P69902, P69903, Q07Q82, Q0TF87, Q139H7, Q1R8Z2, Q217M3, Q2IUI7, Q31Y97, Q32DG9, Q3YZF6, Q6N8F8, Q82M40, Q89QH2, Q8FFE8, Q8XBR7] | R07290: [] | 2.8.3.16_1: [] | 2.8.3.9_14: [] | '
$ws1.Range('G15').Value = 'RXN-863_reverse: [O87942, O87943, O87944, Q8VPT6, Q8VPT8] | R05598: [] | 4.1.99.11_13: [] | '
$ws1.Range('H15').Value = '1.1.3.13_53: [] | '
$ws1.Range('I15').Value = 'RXN0-7141_reverse: [] | '
$ws1.Range('F16').Value = 'RXN-17511_reverse: [] | 1.2.3.1_109: [] | 1.2.3.15_09: [] | 1.17.3.2_40: [] | 1.2.3.1_107: [] | '
$ws1.Range('G16').Value = 'RXN0-1382: [A1ADQ1, A4YXN2, A5EGD7, A7ZPI2, A8A2M8, B1IX88, B1LMH0, B1X9P6, B2TWX3, B3QBS6, B5YYX4, B6I6S5, B6JE29, B7LBS7, B7M6P3, B7MH34, B7MY33, B7N5X4, B7NPQ8, B7UG84, C4ZVR1, O06644, O87838, P69902, P69903, Q07Q82, Q0TF87, Q139H7, Q1R8Z2, Q217M3, Q2IUI7, Q31Y97, Q32DG9, Q3YZF6, Q6N8F8, Q82M40, Q89QH2, Q8FFE8, Q8XBR7] | R07290: [] | 2.8.3.16_1: [] | 2.8.3.9_14: [] | '
$ws1.Range('H16').Value = 'RXN-863_reverse: [O87942, O87943, O87944, Q8VPT6, Q8VPT8] | R05598: [] | 4.1.99.11_13: [] | '
$ws1.Range('I16').Value = 'RXN-16543_reverse: [] | '
$ws1.Range('F17').Value = '2.1.1.15_11: [] | '
$ws1.Range('G17').Value = 'RXN0-1382: [A1ADQ1, A4YXN2, A5EGD7, A7ZPI2, A8A2M8, B1IX88, B1LMH0, B1X9P6, B2TWX3, B3QBS6, B5YYX4, B6I6S5, B6JE29, B7LBS7, B7M6P3, B7MH34, B7MY33, B7N5X4, B7NPQ8, B7UG84, C4ZVR1, O06644, O87838, P69902, P69903, Q07Q82, Q0TF87, Q139H7, Q1R8Z2, Q217M3, Q2IUI7, Q31Y97, Q32DG9, Q3YZF6, Q6N8F8, Q82M40, Q89QH2, Q8FFE8, Q8XBR7] | R07290: [] | 2.8.3.16_1: [] | 2.8.3.9_14: [] | '
$ws1.Range('H17').Value = 'RXN-863_reverse: [O87942, O87943, O87944, Q8VPT6, Q8VPT8] | R05598: [] | 4.1.99.11_13: [] | '
$ws1.Range('I17').Value = 'RXN-16543_reverse: [] | '
$ws1.Range('F18').Value = '1.2.1.4_26: [] | '
$ws1.Range('G18').Value = 'RXN0-1382: [A1ADQ1, A4YXN2, A5EGD7, A7ZPI2, A8A2M8, B1IX88, B1LMH0, B1X9P6, B2TWX3, B3QBS6, B5YYX4, B6I6S5, B6JE29, B7LBS7, B7M6P3, B7MH34, B7MY33, B7N5X4, B7NPQ8, B7UG84, C4ZVR1, O06644, O87838, P69902, P69903, Q07Q82, Q0TF87, Q139H7, Q1R8Z2, Q217M3, Q2IUI7, Q31Y97, Q32DG9, Q3YZF6, Q6N8F8, Q82M40, Q89QH2, Q8FFE8, Q8XBR7] | R07290: [] | 2.8.3.16_1: [] | 2.8.3.9_14: [] | '
$ws1.Range('H18').Value = 'RXN-863_reverse: [O87942, O87943, O87944, Q8VPT6, Q8VPT8] | R05598: [] | 4.1.99.11_13: [] | '
$ws1.Range('I18').Value = 'RXN-16543_reverse: [] | '
$ws1.Range('F19').Value = 'RXN-863_reverse: [O87942, O87943, O87944, Q8VPT6, Q8VPT8] | R05598: [] | 4.1.99.11_13: [] | 4.1.99.11_07: [] | RXN-942_reverse: [Q8VPT6, Q8VPT7] | 4.1.99.11_08: [] | '
$ws1.Range('G19').Value = 'R02243: [] | ACONITATE-DECARBOXYLASE-RXN_reverse: [B3IUN8, P54987, Q0C8L3] | '
$ws1.Range('H19').Value = 'RXN-16543_reverse: [] | '
$ws1.Range('I19').Value = 'R00930_reverse: [] | 2.1.3.1_7: [] | 2.1.3.1-RXN: [P02904, Q70AC7, Q8GBW6] | '
$ws1.Range('F20').Value = 'RXN-863_reverse: [O87942, O87943, O87944, Q8VPT6, Q8VPT8] | R05598: [] | 4.1.99.11_13: [] | 4.1.99.11_07: [] | RXN-942_reverse: [Q8VPT6, Q8VPT7] | 4.1.99.11_08: [] | '
$ws1.Range('G20').Value = 'R02243: [] | ACONITATE-DECARBOXYLASE-RXN_reverse: [B3IUN8, P54987, Q0C8L3] | '
$ws1.Range('H20').Value = 'RXN-16543_reverse: [] | '
$ws1.Range('I20').Value = 'RXN0-1382: [A1ADQ1, A4YXN2, A5EGD7, A7ZPI2, A8A2M8, B1IX88, B1LMH0, B1X9P6, B2TWX3, B3QBS6, B5YYX4, B6I6S5, B6JE29, B7LBS7, B7M6P3, B7MH34, B7MY33, B7N5X4, B7NPQ8, B7UG84, C4ZVR1, O06644, O87838, P69902, P69903, Q07Q82, Q0TF87, Q139H7, Q1R8Z2, Q217M3, Q2IUI7, Q31Y97, Q32DG9, Q3YZF6, Q6N8F8, Q82M40, Q89QH2, Q8FFE8, Q8XBR7] | R07290: [] | 2.8.3.16_1: [] | '
$ws1.Range('F21').Value = 'R01776: [] | HOMOSERINE-O-ACETYLTRANSFERASE-RXN_reverse: [A0A0D2YG05, A0KBH0, A0LCI7, A0PNQ1, A0QKL5, A0QSZ0, A1AUS2, A1KCN1, A1KP32, A1SMU5, A1TY14, A1UC78, A1V7J8, A1WVL3, A2S877, A2SLU8, A3M1Y2, A3MRS1, A3N4G4, A3NQ54, A4G1K0, A4JIS1, A4VRJ2, A4XZY6, A5CQM6, A5FVN0, A5GC94, A5U815, A5WAD0, A6SUC4, A6UYJ9, A6W3E8, A9AC42, A9I0E6, B0KN51, B0RS70, B0SHJ0, B0SQZ8, B0U3V7, B0V4S0, B0VL51, B1J2H5, B1K0G8, B1VF43, B1XW69, B1YPU9, B2AGC8, B2FMC1, B2HDS6, B2I3B6, B2I700, B2IY96, B2JJW9, B2T1N1, B2U7Y8, B3E278, B3PFK3, B3QM89, B3QZ51, B4E7Y5, B5EAS1, B5ELV1, B7H098, B7I4X4, B7J437, B7V3X6, B8GZZ0, B8ZUV9, B9M4M5, C1AHF8, C1D4Z8, C1DI67, C3K3L0, C5BLJ8, O27848, O32874, O60062, O68640, P08465, P0A5J9, P12917, P45131, P57714, P57715, P57884, P94891, P9WJY8, P9WJY9, Q02U23, Q04RK8, Q052B3, Q06736, Q0ABU4, Q0AIM3, Q0BAW3, Q0BVW2, Q0KF58, Q146X7, Q1BCR5, Q1BSM1, Q1IGC6, Q1LS47, Q1QQD3, Q20YI1, Q21EI0, Q2IS70, Q2KU63, Q2P3Z2, Q2RIW8, Q2RP84, Q2S8M0, Q2T284, Q2YCJ2, Q31EJ2, Q39BX2, Q39RX3, Q3A7G2, Q3AAW2, Q3ARA7, Q3B5A5, Q3BSK1, Q3IUE8, Q3JXB2, Q3K594, Q3M5Q6, Q3SV39, Q476W5, Q47JU2, Q47L17, Q48PA1, Q4FP51, Q4FUR4, Q4JTI5, Q4K4D3, Q4QJY8, Q4UVH1, Q4ZZ78, Q5F858, Q5FUF4, Q5H124, Q5NR05, Q5P4R1, Q5SK89, Q5UY71, Q5Z1C4, Q62F07, Q63YJ0, Q65U29, Q67NS3, Q6ADB9, Q6ANV2, Q6FEQ3, Q6N1G8, Q6NIZ3, Q722V6, Q72R95, Q73UB0, Q74AC8, Q7CWE8, Q7NZY3, Q7W0P6, Q7W3F8, Q7WES8, Q81ZZ5, Q87BG9, Q87V90, Q88CT3, Q89UL7, Q8F4I0, Q8FRT0, Q8KES7, Q8P8L2, Q8PK44, Q8Y3F3, Q8Y9D6, Q92E58, Q98G09, Q9AAS1, Q9JUT9, Q9JZQ5, Q9PAN0, Q9RA51, Q9RVZ8, Q9Y875, S0DUX2, W7N293] | 2.3.1.31_01: [] | R01023: [] | CHOLINE-O-ACETYLTRANSFERASE-RXN_reverse: [B2ZGJ1, P07668, P13222, P28329, P32738, P32756, Q03059, Q90YJ9] | 2.3.1.7_19: [] | 2.3.1.6_2: [] | 2.3.1.84_17: [] | 2.3.1.84_10: [] | '
$ws1.Range('G21').Value = 'RXN0-1382: [A1ADQ1, A4YXN2, A5EGD7, A7ZPI2, A8A2M8, B1IX88, B1LMH0, B1X9P6, B2TWX3, B3QBS6, B5YYX4, B6I6S5, B6JE29, B7LBS7, B7M6P3, B7MH34, B7MY33, B7N5X4, B7NPQ8, B7UG84, C4ZVR1, O06644, O87838, P69902, P69903, Q07Q82, Q0TF87, Q139H7, Q1R8Z2, Q217M3, Q2IUI7, Q31Y97, Q32DG9, Q3YZF6, Q6N8F8, Q82M40, Q89QH2, Q8FFE8, Q8XBR7] | R07290: [] | 2.8.3.16_1: [] | 2.8.3.9_14: [] | '
$ws1.Range('H21').Value = 'RXN-863_reverse: [O87942, O87943, O87944, Q8VPT6, Q8VPT8] | R05598: [] | 4.1.99.11_13: [] | '
$ws1.Range('I21').Value = 'RXN-16543_reverse: [] | '
$ws1.Range('F22').Value = 'R01776: [] | HOMOSERINE-O-ACETYLTRANSFERASE-RXN_reverse: [A0A0D2YG05, A0KBH0, A0LCI7, A0PNQ1, A0QKL5, A0QSZ0, A1AUS2, A1KCN1, A1KP32, A1SMU5, A1TY14, A1UC78, A1V7J8, A1WVL3, A2S877, A2SLU8, A3M1Y2, A3MRS1, A3N4G4, A3NQ54, A4G1K0, A4JIS1, A4VRJ2, A4XZY6, A5CQM6, A5FVN0, A5GC94, A5U815, A5WAD0, A6SUC4, A6UYJ9, A6W3E8, A9AC42, A9I0E6, B0KN51, B0RS70, B0SHJ0, B0SQZ8, B0U3V7, B0V4S0, B0VL51, B1J2H5, B1K0G8, B1VF43, B1XW69, B1YPU9, B2AGC8, B2FMC1, B2HDS6, B2I3B6, B2I700, B2IY96, B2JJW9, B2T1N1, B2U7Y8, B3E278, B3PFK3, B3QM89, B3QZ51, B4E7Y5, B5EAS1, B5ELV1, B7H098, B7I4X4, B7J437, B7V3X6, B8GZZ0, B8ZUV9, B9M4M5, C1AHF8, C1D4Z8, C1DI67, C3K3L0, C5BLJ8, O27848, O32874, O60062, O68640, P08465, P0A5J9, P12917, P45131, P57714, P57715, P57884, P94891, P9WJY8, P9WJY9, Q02U23, Q04RK8, Q052B3, Q06736, Q0ABU4, Q0AIM3, Q0BAW3, Q0BVW2, Q0KF58, Q146X7, Q1BCR5, Q1BSM1, Q1IGC6, Q1LS47, Q1QQD3, Q20YI1, Q21EI0, Q2IS70, Q2KU63, Q2P3Z2, Q2RIW8, Q2RP84, Q2S8M0, Q2T284, Q2YCJ2, Q31EJ2, Q39BX2, Q39RX3, Q3A7G2, Q3AAW2, Q3ARA7, Q3B5A5, Q3BSK1, Q3IUE8, Q3JXB2, Q3K594, Q3M5Q6, Q3SV39, Q476W5, Q47JU2, Q47L17, Q48PA1, Q4FP51, Q4FUR4, Q4JTI5, Q4K4D3, Q4QJY8, Q4UVH1, Q4ZZ78, Q5F858, Q5FUF4, Q5H124, Q5NR05, Q5P4R1, Q5SK89, Q5UY71, Q5Z1C4, Q62F07, Q63YJ0, Q65U29, Q67NS3, Q6ADB9, Q6ANV2, Q6FEQ3, Q6N1G8, Q6NIZ3, Q722V6, Q72R95, Q73UB0, Q74AC8, Q7CWE8, Q7NZY3, Q7W0P6, Q7W3F8, Q7WES8, Q81ZZ5, Q87BG9, Q87V90, Q88CT3, Q89UL7, Q8F4I0, Q8FRT0, Q8KES7, Q8P8L2, Q8PK44, Q8Y3F3, Q8Y9D6, Q92E58, Q98G09, Q9AAS1, Q9JUT9, Q9JZQ5, Q9PAN0, Q9RA51, Q9RVZ8, Q9Y875, S0DUX2, W7N293] | 2.3.1.31_01: [] | R01023: [] | CHOLINE-O-ACETYLTRANSFERASE-RXN_reverse: [B2ZGJ1, P07668, P13222, P28329, P32738, P32756, Q03059, Q90YJ9] | 2.3.1.7_19: [] | 2.3.1.6_2: [] | 2.3.1.84_17: [] | 2.3.1.84_10: [] | '
$ws1.Range('G22').Value = 'RXN0-1382: [A1ADQ1, A4YXN2, A5EGD7, A7ZPI2, A8A2M8, B1IX88, B1LMH0, B1X9P6, B2TWX3, B3QBS6, B5YYX4, B6I6S5, B6JE29, B7LBS7, B7M6P3, B7MH34, B7MY33, B7N5X4, B7NPQ8, B7UG84, C4ZVR1, O06644, O87838, P69902, P69903, Q07Q82, Q0TF87, Q139H7, Q1R8Z2, Q217M3, Q2IUI7, Q31Y97, Q32DG9, Q3YZF6, Q6N8F8, Q82M40, Q89QH2, Q8FFE8, Q8XBR7] | R07290: [] | 2.8.3.16_1: [] | 2.8.3.9_14: [] | '
$ws1.Range('H22').Value = 'RXN-863_reverse: [O87942, O87943, O87944, Q8VPT6, Q8VPT8] | R05598: [] | 4.1.99.11_13: [] | '
$ws1.Range('I22').Value = 'RXN-16543_reverse: [] | '
$ws1.Range('F23').Value = 'RXN-11989: [] | R09782_reverse: [] | 1.13.11.67_0: [] | 1.13.11.59_1: [] | '
$ws1.Range('G23').Value = 'RXN0-1382: [A1ADQ1, A4YXN2, A5EGD7, A7ZPI2, A8A2M8, B1IX88, B1LMH0, B1X9P6, B2TWX3, B3QBS6, B5YYX4, B6I6S5, B6JE29, B7LBS7, B7M6P3, B7MH34, B7MY33, B7N5X4, B7NPQ8, B7UG84, C4ZVR1, O06644, O87838, P69902, P69903, Q07Q82, Q0TF87, Q139H7, Q1R8Z2, Q217M3, Q2IUI7, Q31Y97, Q32DG9, Q3YZF6, Q6N8F8, Q82M40, Q89QH2, Q8FFE8, Q8XBR7] | R07290: [] | 2.8.3.16_1: [] | 2.8.3.9_14: [] | '
$ws1.Range('H23').Value = 'RXN-863_reverse: [O87942, O87943, O87944, Q8VPT6, Q8VPT8] | R05598: [] | 4.1.99.11_13: [] | '
$ws1.Range('I23').Value = 'RXN-16543_reverse: [] | '
$ws1.Range('F24').Value = 'R01440_reverse: [] | FORMALDEHYDE-TRANSKETOLASE-RXN: [P06834] | 2.2.1.3_16: [] | 2.2.1.3_15: [] | '
$ws1.Range('G24').Value = 'RXN0-1382: [A1ADQ1, A4YXN2, A5EGD7, A7ZPI2, A8A2M8, B1IX88, B1LMH0, B1X9P6, B2TWX3, B3QBS6, B5YYX4, B6I6S5, B6JE29, B7LBS7, B7M6P3, B7MH34, B7MY33, B7N5X4, B7NPQ8, B7UG84, C4ZVR1, O06644, O87838, P69902, P69903, Q07Q82, Q0TF87, Q139H7, Q1R8Z2, Q217M3, Q2IUI7, Q31Y97, Q32DG9, Q3YZF6, Q6N8F8, Q82M40, Q89QH2, Q8FFE8, Q8XBR7] | R07290: [] | 2.8.3.16_1: [] | 2.8.3.9_14: [] | '
$ws1.Range('H24').Value = 'RXN-863_reverse: [O87942, O87943, O87944, Q8VPT6, Q8VPT8] | R05598: [] | 4.1.99.11_13: [] | '
$ws1.Range('I24').Value = 'RXN-16543_reverse: [] | '

$ws2 = $wb.Worksheets.Item("pyruvate_mvacid")
$ws2.Range('F2').Value = '1.4.3.3_11: [] | 1.4.3.2_22: [] | '
$ws2.Range('G2').Value = '5.1.1.10_44: [] | 4.3.1.17_15: [] | '
$ws2.Range('H2').Value = 'RXN-863_reverse: [O87942, O87943, O87944, Q8VPT6, Q8VPT8] | R05598: [] | 4.1.99.11_13: [] | '
$ws2.Range('I2').Value = 'RXN-16543_reverse: [] | '
$ws2.Range('F3').Value = '2.6.1.42_063: [] | 2.6.1.42_064: [] | '
$ws2.Range('G3').Value = '5.1.1.10_44: [] | 4.3.1.17_15: [] | '
$ws2.Range('H3').Value = 'RXN-863_reverse: [O87942, O87943, O87944, Q8VPT6, Q8VPT8] | R05598: [] | 4.1.99.11_13: [] | '
$ws2.Range('I3').Value = 'RXN-16543_reverse: [] | '
